# Update the "Journal" work-log sheet:
#  - Bulle/Domaine/Module header (F3) changes topic
#  - Row 6 gets corrected duration / status / end-time
#  - Rows 7 and 8 (previously blank placeholder rows) get filled in with two
#    new work-log entries ("Class diagram" and "Explanations")
#  - Totals (D11, C54) are plain formulas and recalculate on their own.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# --- Header block -----------------------------------------------------
# F3:G3 — "Bulle / Domaine / Module" value
$ws.Range("F3").Value = "INF / DEV / P_Prod"

# --- Row 6 : existing entry gets corrected ----------------------------
$ws.Range("C6").Value = 30
$ws.Range("E6").Value = "Finished"
$ws.Range("F6").Value = 0.59027777777777779

# --- Row 8 : new entry "Class diagram" --------------------------------
$ws.Range("A8").Value = "Class diagram"
$ws.Range("C8").Value = 45
$ws.Range("D8").Value = "Creating a diagram for the code's classes"
$ws.Range("E8").Value = "In the work"

$ws.Range("F8").Value = 0.65625
$ws.Range("F8").NumberFormat = "h:mm"

# --- Row 7 : new entry "Explanations" ---------------------------------
$ws.Range("A7").Value = "Explanations"
$ws.Range("A7").WrapText = $false
$ws.Range("A7").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A7").Borders.LineStyle = -4142     # xlLineStyleNone

$ws.Range("C7").Value = 35

$ws.Range("D7").Value = "Explanations of the code for me"
$ws.Range("D7").WrapText = $false
$ws.Range("D7").HorizontalAlignment = 1       # xlGeneral
$ws.Range("D7").Borders.LineStyle = -4142     # xlLineStyleNone

$ws.Range("E7").Value = "Finished"

$ws.Range("F7").Value = 0.61458333333333337
$ws.Range("F7").NumberFormat = "h:mm"

# G8 link must be the very last new shared string created.
$ws.Range("G8").Value = "https://eduvaud-my.sharepoint.com/:u:/r/personal/pb62kjx_eduvaud_ch/_layouts/15/Doc.aspx?sourcedoc=%7B5BF0A89E-E3A6-4CD4-BA7E-BFDC17A68BD3%7D&file=P_CraftMeUp_ClassDiagram.vsdx&fromShare=true&action=default&mobileredirect=true"

# --- Recalculate totals & refresh the remembered selection ------------
$wb.Application.Calculate()
$ws.Range("I8").Select()
